# edit.ps1 - Applies the cryptos.xlsx data refresh described by the commit diff.
# Updates coin prices (column D) and 1h volume change percentages (column E)
# for rows 2-51, plus a name/link/price/volume swap for rows 32/33 (Cosmos<->Hedera)
# and rows 41/42 (Stacks<->Kaspa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $val into $cellRef while forcing a *text* result, even when
# $val looks like a number (e.g. "9.70", "0.120"). A direct `.Value =` assignment
# would let Excel coerce such strings into real numbers (losing the formatting/
# trailing zeros and changing the stored type), so instead we push the text in
# via a `="..."` formula and immediately collapse it to a plain value with a
# copy / paste-values, which keeps the cell a plain string cell.
function Set-TextValue($cellRef, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
}

# Row 2
$ws.Range("D2").Value = '70.257.12'
$ws.Range("E2").Value = '  -0.10%  '
# Row 3
$ws.Range("D3").Value = '3.604.99'
$ws.Range("E3").Value = '  -0.18%  '
# Row 4
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.05%  '
# Row 5
Set-TextValue "D5" '580.72'
$ws.Range("E5").Value = '  -1.32%  '
# Row 6
Set-TextValue "D6" '192.26'
$ws.Range("E6").Value = '  +0.69%  '
# Row 7
Set-TextValue "D7" '0.633'
$ws.Range("E7").Value = '  -1.61%  '
# Row 8
$ws.Range("D8").Value = '3.602.39'
$ws.Range("E8").Value = '  -0.02%  '
# Row 9
$ws.Range("E9").Value = '  +0.01%  '
# Row 10
$ws.Range("E10").Value = '  +2.68%  '
# Row 11
Set-TextValue "D11" '0.667'
$ws.Range("E11").Value = '  +0.95%  '
# Row 12
Set-TextValue "D12" '56.08'
$ws.Range("E12").Value = '  -3.61%  '
# Row 13
Set-TextValue "D13" '0.0000309'
$ws.Range("E13").Value = '  +6.08%  '
# Row 14
Set-TextValue "D14" '9.70'
$ws.Range("E14").Value = '  -0.99%  '
# Row 15
$ws.Range("D15").Value = '4.173.66'
$ws.Range("E15").Value = '  -0.33%  '
# Row 16
Set-TextValue "D16" '19.96'
$ws.Range("E16").Value = '  +3.02%  '
# Row 17
$ws.Range("D17").Value = '3.593.90'
$ws.Range("E17").Value = '  -0.37%  '
# Row 18
$ws.Range("D18").Value = '70.210.72'
$ws.Range("E18").Value = '  +0.02%  '
# Row 19
Set-TextValue "D19" '12.71'
$ws.Range("E19").Value = '  +1.75%  '
# Row 20
$ws.Range("E20").Value = '  +0.31%  '
# Row 21
$ws.Range("E21").Value = '  -0.31%  '
# Row 22
Set-TextValue "D22" '480.29'
$ws.Range("E22").Value = '  -2.51%  '
# Row 23
Set-TextValue "D23" '19.21'
$ws.Range("E23").Value = '  +10.77%  '
# Row 24
$ws.Range("E24").Value = '  -6.50%  '
# Row 25
Set-TextValue "D25" '4.39'
$ws.Range("E25").Value = '  -1.35%  '
# Row 26
Set-TextValue "D26" '95.47'
$ws.Range("E26").Value = '  +5.24%  '
# Row 27
$ws.Range("E27").Value = '  -3.34%  '
# Row 28
Set-TextValue "D28" '11.07'
$ws.Range("E28").Value = '  -0.46%  '
# Row 29
Set-TextValue "D29" '9.34'
$ws.Range("E29").Value = '  -1.32%  '
# Row 30
Set-TextValue "D30" '32.34'
$ws.Range("E30").Value = '  -0.26%  '
# Row 31
Set-TextValue "D31" '7.71'
$ws.Range("E31").Value = '  +2.19%  '
# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D32" '0.120'
$ws.Range("E32").Value = '  +2.23%  '
# Row 33
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D33" '12.25'
$ws.Range("E33").Value = '  -0.01%  '
# Row 34
Set-TextValue "D34" '66.67'
$ws.Range("E34").Value = '  +1.97%  '
# Row 35
Set-TextValue "D35" '588.60'
$ws.Range("E35").Value = '  -6.14%  '
# Row 36
Set-TextValue "D36" '39.30'
$ws.Range("E36").Value = '  +3.15%  '
# Row 38
$ws.Range("E38").Value = '  -1.82%  '
# Row 39
Set-TextValue "D39" '0.397'
$ws.Range("E39").Value = '  -1.84%  '
# Row 40
Set-TextValue "D40" '3.24'
$ws.Range("E40").Value = '  +19.54%  '
# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D41" '0.138'
$ws.Range("E41").Value = '  -5.39%  '
# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D42" '3.47'
$ws.Range("E42").Value = '  -4.05%  '
# Row 43
$ws.Range("D43").Value = '3.239.01'
$ws.Range("E43").Value = '  -1.87%  '
# Row 44
Set-TextValue "D44" '2.87'
$ws.Range("E44").Value = '  +6.94%  '
# Row 45
Set-TextValue "D45" '3.08'
$ws.Range("E45").Value = '  -0.43%  '
# Row 46
Set-TextValue "D46" '0.0448'
$ws.Range("E46").Value = '  +0.60%  '
# Row 47
Set-TextValue "D47" '9.50'
$ws.Range("E47").Value = '  +4.62%  '
# Row 48
Set-TextValue "D48" '3.34'
$ws.Range("E48").Value = '  -0.54%  '
# Row 49
$ws.Range("E49").Value = '  +0.60%  '
# Row 50
$ws.Range("E50").Value = '  -0.13%  '
# Row 51
Set-TextValue "D51" '3.16'
$ws.Range("E51").Value = '  -4.38%  '
